$d = $word.ActiveDocument

# The paragraph describing the dataset features originally has the F1..F11
# feature labels inline as plain text. The edit bolds each "Fn - <type>"
# label (splitting the existing runs) while leaving the surrounding
# descriptive text unchanged. No characters are inserted/removed - this is
# purely a run-splitting / character-formatting change.

$targets = @(
    "F1 - Continuous value",
    "F2 - Continuous value",
    "F4 - Categorical Value",
    "F5 - continuous value",
    "F6 - continuous value",
    "F7 - Categorical value",
    "F8 - Categorical value",
    "F9 Categorical Value",
    "F10 - Categorical Value",
    "F11 - Categorical"
)

# Start the search right at the beginning of the paragraph that contains
# the feature list, so we don't accidentally match similar text elsewhere
# in the document.
$anchor = $d.Content
$anchor.Find.Execute("The features are, in order", $true, $false, $false, $false, $false,
                      $true, 1, $false, "", 0) | Out-Null

$cursor = $d.Range($anchor.Start, $anchor.Start)

foreach ($target in $targets) {
    $found = $cursor.Find.Execute($target, $true, $false, $false, $false, $false,
                                   $true, 1, $false, "", 0)
    if ($found) {
        $cursor.Bold = 1
        $cursor.BoldBi = 1
        $cursor = $d.Range($cursor.End, $cursor.End)
    }
}
